$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.234.29'
$ws.Range("E2").Value = '  -3.70%  '

$ws.Range("D3").Value = '2.463.00'
$ws.Range("E3").Value = '  -2.99%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.32'
$ws.Range("E5").Value = '  +0.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.21'
$ws.Range("E6").Value = '  -6.60%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.547'
$ws.Range("E7").Value = '  -3.42%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("E9").Value = '  -4.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.42'
$ws.Range("E10").Value = '  -6.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0781'
$ws.Range("E11").Value = '  -3.03%  '

$ws.Range("E12").Value = '  -0.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.98'
$ws.Range("E13").Value = '  -4.84%  '

$ws.Range("D14").Value = '2.841.26'
$ws.Range("E14").Value = '  -3.05%  '

$ws.Range("D15").Value = '2.482.82'
$ws.Range("E15").Value = '  -2.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.96'
$ws.Range("E16").Value = '  -2.47%  '

$ws.Range("E17").Value = '  -3.73%  '

$ws.Range("D18").Value = '41.228.00'
$ws.Range("E18").Value = '  -3.71%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.31'
$ws.Range("E19").Value = '  -6.53%  '

$ws.Range("D20").Value = '0.0₃0922'
$ws.Range("E20").Value = '  -3.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.22'
$ws.Range("E21").Value = '  -9.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.35'
$ws.Range("E22").Value = '  -2.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.50'
$ws.Range("E23").Value = '  -2.55%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.74'
$ws.Range("E24").Value = '  -5.07%  '

$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("E26").Value = '  -6.72%  '

$ws.Range("E27").Value = '  -6.14%  '

$ws.Range("E28").Value = '  -4.31%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.66'
$ws.Range("E29").Value = '  -5.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.37'
$ws.Range("E30").Value = '  -5.84%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '151.54'
$ws.Range("E31").Value = '  -4.50%  '

$ws.Range("E32").Value = '  -6.76%  '

$ws.Range("E33").Value = '  -5.84%  '

$ws.Range("E34").Value = '  -3.21%  '

$ws.Range("E35").Value = '  -6.18%  '

$ws.Range("E36").Value = '  -3.26%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.30'
$ws.Range("E37").Value = '  -4.79%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.86'
$ws.Range("E38").Value = '  -5.54%  '

$ws.Range("E39").Value = '  -3.03%  '

$ws.Range("E40").Value = '  -9.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.22'
$ws.Range("E41").Value = '  +1.75%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.52'
$ws.Range("E43").Value = '  -10.84%  '

$ws.Range("D44").Value = '1.982.85'
$ws.Range("E44").Value = '  -0.79%  '

$ws.Range("E45").Value = '  -4.86%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.01'
$ws.Range("E46").Value = '  -9.20%  '

$ws.Range("E47").Value = '  -5.11%  '

$ws.Range("D48").Value = '2.704.88'
$ws.Range("E48").Value = '  -2.81%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '69.54'
$ws.Range("E49").Value = '  -4.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '96.52'
$ws.Range("E50").Value = '  -4.70%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.45'
$ws.Range("E51").Value = '  -7.25%  '

